$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update footer timestamp (row 1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 21 de Marzo de 2020 a las 01:16"

# Update reordered country names and refreshed case numbers
$ws.Cells.Item(4, 2).Value = 81008
$ws.Cells.Item(4, 3).Value = 41
$ws.Cells.Item(4, 4).Value = 71740
$ws.Cells.Item(4, 5).Value = 6013
$ws.Cells.Item(4, 6).Value = 1927
$ws.Cells.Item(4, 7).Value = 7
$ws.Cells.Item(4, 8).Value = 3255
$ws.Cells.Item(9, 2).Value = 19393
$ws.Cells.Item(9, 3).Value = 5604
$ws.Cells.Item(9, 5).Value = 18990
$ws.Cells.Item(17, 2).Value = 1957
$ws.Cells.Item(17, 3).Value = 167
$ws.Cells.Item(17, 5).Value = 1949
$ws.Cells.Item(106, 1).Value = "Liechtenstein"
$ws.Cells.Item(107, 1).Value = "Reunion"
$ws.Cells.Item(122, 1).Value = "Montenegro"
$ws.Cells.Item(122, 3).Value = 1
$ws.Cells.Item(123, 1).Value = "Guam"
$ws.Cells.Item(123, 3).Value = 2
$ws.Cells.Item(129, 1).Value = "Monaco"
$ws.Cells.Item(129, 3).Value = 1
$ws.Cells.Item(130, 1).Value = "Polinesia Francesa"
$ws.Cells.Item(130, 3).Value = 5
$ws.Cells.Item(133, 1).Value = "Etiopia"
$ws.Cells.Item(133, 3).Value = 2
$ws.Cells.Item(134, 1).Value = "Togo"
$ws.Cells.Item(134, 3).Value = 8
$ws.Cells.Item(136, 1).Value = "Kenia"
$ws.Cells.Item(136, 3).Value = 0
$ws.Cells.Item(137, 1).Value = "Seychelles"
$ws.Cells.Item(137, 3).Value = 1
$ws.Cells.Item(138, 1).Value = "Kirguistan"
$ws.Cells.Item(138, 3).Value = 3
$ws.Cells.Item(139, 1).Value = "Mayotte"
$ws.Cells.Item(139, 3).Value = 2
$ws.Cells.Item(140, 1).Value = "Tanzania"
$ws.Cells.Item(140, 3).Value = 0
$ws.Cells.Item(141, 1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(142, 1).Value = "Barbados"
$ws.Cells.Item(142, 3).Value = 1
$ws.Cells.Item(149, 1).Value = "San Bartolome"
$ws.Cells.Item(150, 1).Value = "Congo"
$ws.Cells.Item(151, 1).Value = "Namibia"
$ws.Cells.Item(151, 3).Value = 0
$ws.Cells.Item(152, 1).Value = "Islas Virgenes de los Estados Unidos"
$ws.Cells.Item(154, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(154, 3).Value = 2
$ws.Cells.Item(155, 1).Value = "Curazao"
$ws.Cells.Item(156, 1).Value = "Islas Caimanes"
$ws.Cells.Item(157, 1).Value = "Nicaragua"
$ws.Cells.Item(158, 1).Value = "Benin"
$ws.Cells.Item(159, 1).Value = "Liberia"
$ws.Cells.Item(159, 3).Value = 0
$ws.Cells.Item(160, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(160, 3).Value = 0
$ws.Cells.Item(162, 1).Value = "Zambia"
$ws.Cells.Item(163, 1).Value = "Mauritania"
$ws.Cells.Item(164, 1).Value = "Guinea"
$ws.Cells.Item(164, 3).Value = 1
$ws.Cells.Item(165, 1).Value = "Butan"
$ws.Cells.Item(165, 3).Value = 1
$ws.Cells.Item(166, 1).Value = "Haiti"
$ws.Cells.Item(166, 3).Value = 2
$ws.Cells.Item(167, 1).Value = "Groenlandia"
$ws.Cells.Item(167, 3).Value = 0
$ws.Cells.Item(168, 1).Value = "Isla de Man"
$ws.Cells.Item(168, 3).Value = 1
$ws.Cells.Item(169, 1).Value = "Bermudas"
$ws.Cells.Item(169, 3).Value = 0
$ws.Cells.Item(171, 1).Value = "Gambia"
$ws.Cells.Item(171, 3).Value = 0
$ws.Cells.Item(172, 1).Value = "Montserrat"
$ws.Cells.Item(173, 1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(174, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(175, 1).Value = "El Salvador"
$ws.Cells.Item(176, 1).Value = "Suazilandia"
$ws.Cells.Item(177, 1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(178, 1).Value = "Fiyi"
$ws.Cells.Item(179, 1).Value = "Republica del Chad"
$ws.Cells.Item(180, 1).Value = "Niger"
$ws.Cells.Item(181, 1).Value = "Santa Sede"
$ws.Cells.Item(182, 1).Value = "Somalia"
$ws.Cells.Item(183, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(184, 1).Value = "Cabo Verde"
$ws.Cells.Item(184, 3).Value = 1
$ws.Cells.Item(185, 1).Value = "Angola"
$ws.Cells.Item(186, 1).Value = "Zimbabue"

Write-Host "Update complete"
